$wb = $excel.ActiveWorkbook

# --- shulConfigeration sheet ---
$ws4 = $wb.Worksheets.Item("shulConfigeration")
$ws4.Activate()

# Cells D89:D91 currently hold plain numbers (905, 915, 925) styled like the
# other formula-result cells in column D. The edit re-formats them as text
# (number format "@" / built-in numFmtId 49) with centered horizontal
# alignment, matching the look of the surrounding text cells, without
# altering the stored numeric value ("change numbers to strings" look).
$rngD89_91 = $ws4.Range("D89:D91")
$rngD89_91.NumberFormat = "@"
$rngD89_91.HorizontalAlignment = -4108

# Scroll/selection moved from J103 to F113.
$ws4.Range("F113").Select()
$win4 = $excel.ActiveWindow
$win4.ScrollRow = 75
$win4.ScrollColumn = 1

# --- shulConfigeration2016 sheet ---
$ws7 = $wb.Worksheets.Item("shulConfigeration2016")
$ws7.Activate()

# Selection/active cell (G102) is unchanged; only the viewport scrolled.
$ws7.Range("G102").Select()
$win7 = $excel.ActiveWindow
$win7.ScrollRow = 68
$win7.ScrollColumn = 1

# Restore shulConfigeration as the active/selected sheet and tab.
$ws4.Activate()
